# findNachbarn und zaehle beide zusammengefasst in einem Code und Methoden Aufruf entfernt
#
# Adds four new log rows (rows 8-11) to the "Zeiten" tracking sheet, mirroring
# the existing Janes/Elias time-tracking entries in columns A (author), B
# (task description) and C (hours).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Elias - printSpielfeld vereinfacht - 24.3h
$ws.Cells.Item(8, 1).Value = "Elias"
$ws.Cells.Item(8, 2).Value = "printSpielfeld vereinfacht"
$ws.Cells.Item(8, 3).Value = 24.3

# Row 10 - Elias - Compiler Optimierung cmd - 14.9h
$ws.Cells.Item(10, 1).Value = "Elias"
$ws.Cells.Item(10, 2).Value = "Compiler Optimierung cmd"
$ws.Cells.Item(10, 3).Value = 14.9

# Row 11 - Janes - findNachbarn u. zaehleLebnde zusammengefasst u. Methoden Aufruf entfernt - 7.7h
$ws.Cells.Item(11, 1).Value = "Janes"
$ws.Cells.Item(11, 2).Value = "findNachbarn u. zaehleLebnde zusammengefasst u. Methoden Aufruf entfernt"
$ws.Cells.Item(11, 3).Value = 7.7

# Row 9 - Elias - pruefeRegeln vereinfacht u. Methoden Aufruf entfernt - 23.9h
$ws.Cells.Item(9, 1).Value = "Elias"
$ws.Cells.Item(9, 2).Value = "pruefeRegeln vereinfacht u. Methoden Aufruf entfernt"
$ws.Cells.Item(9, 3).Value = 23.9

# Column B now holds much longer descriptions - widen it (AutoFit / best-fit
# resize) to fit the new text, matching what Excel does when you double-click
# the column border after typing longer values.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 65.109375

# Move/restore the active selection like Excel would after entering the data.
[void]$ws.Range("B16").Select()
